$d = $word.ActiveDocument

# The "Experience" section has the "Public Broadcasting Service" Heading2
# duplicated immediately before both the "Manager of Digital Development"
# and "Senior Web Technologist" Heading3 entries (it should only appear
# once, above "Senior Manager of Technology Solutions"). Find every
# "Public Broadcasting Service" Heading 2 paragraph and drop every
# occurrence after the first one.
$targets = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.Trim() -eq "Public Broadcasting Service") {
        $targets += $i
    }
}

for ($j = $targets.Count - 1; $j -ge 1; $j--) {
    $d.Paragraphs.Item($targets[$j]).Range.Delete()
}

# Fix a stray ", " left over before "ColdFusion" in the Senior Web
# Technologist tech-stack line.
$d.Content.Find.Execute("2005 - 2007 " + [char]0x2014 + " , ColdFusion", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2005 - 2007 " + [char]0x2014 + " ColdFusion", 2)
